$wb = $excel.ActiveWorkbook

# --- Update template version on the "isa_template" metadata sheet ---
$meta = $wb.Worksheets.Item("isa_template")
$meta.Range("B4").Value = "1.0.2"

# --- Add example values to the annotation table on "New Table" sheet ---
$data = $wb.Worksheets.Item("New Table")

$data.Range("B2").Value = "Gas Chromatography"
$data.Range("C2").Value = "NCIT"
$data.Range("D2").Value = "http://purl.obolibrary.org/obo/NCIT_C30014"
$data.Range("F2").Value = "2 ml ethyl acetate"
$data.Range("I2").Value = "trimethylsilyl derivatisation"
$data.Range("J2").Value = "CHMO"
$data.Range("K2").Value = "http://purl.obolibrary.org/obo/CHMO_0002758"
$data.Range("L2").Value = "Shimadzu GCMS-QP2010 Ultra"
$data.Range("O2").Value = "Zebron ZB-AAA GC(10 m x 0.25 mm; Phenomenex)"
$data.Range("R2").Value = "medium polarity"
